$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.248.89'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.606.72'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.44'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  +0.44%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.23'
$ws.Range("E10").Value = '  -0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0810'
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.828.61'
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.596.59'
$ws.Range("E13").Value = '  -0.99%  '
$ws.Range("E14").Value = '  -0.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.517'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.231.52'
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.42'
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.80'
$ws.Range("E20").Value = '  +2.20%  '
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.28'
$ws.Range("E22").Value = '  -2.28%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.95'
$ws.Range("E24").Value = '  +11.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.59'
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  -7.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.22'
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.56'
$ws.Range("E29").Value = '  +0.74%  '
$ws.Range("E30").Value = '  +4.00%  '
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("E34").Value = '  -1.95%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.138.25'
$ws.Range("E36").Value = '  +2.65%  '
$ws.Range("E37").Value = '  +6.16%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("E41").Value = '  -1.93%  '
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("E43").Value = '  +0.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.741.47'
$ws.Range("E44").Value = '  -0.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.07'
$ws.Range("E45").Value = '  -1.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.50'
$ws.Range("E46").Value = '  -3.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.23'
$ws.Range("E48").Value = '  -0.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.408'
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("E51").Value = '  -9.39%  '
